# Daily attendance processing - 2025-12-31 05:11:07
# Reorders the "Recorded By" (column G) contributor list for each session row:
# the comma-separated list of recorders is reversed in order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($null -eq $val) { continue }

    $parts = $val -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -gt 1) {
        $reversed = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = [string]::Join(', ', $reversed)
    }
}
